# Finance_Project_Charter.docx content rebrand: AI/ML Implementation -> Finance/Banking
# Core Banking System Modernization, plus enabling a (blank) default header/footer.

$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title block
Replace-All "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING" "FINANCE - CORE BANKING SYSTEM MODERNIZATION"

# Subtitle (also reused verbatim later in the Project Overview section)
Replace-All "Strategic Initiative for Digital transformation through intelligent automation and predictive analytics" "Strategic Initiative for Digital transformation through modern banking operations and transaction processing"

# Document info block
Replace-All "Industry: Finance and Machine Learning" "Industry: Banking and Banking Operations"
Replace-All "Project Type: Finance Implementation" "Project Type: Banking Implementation"

# Executive summary paragraph
Replace-All "This project proposal outlines a strategic Finance Implementation initiative for Finance and Machine Learning to achieve Digital transformation through intelligent automation and predictive analytics. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization." "This project proposal outlines a strategic Banking Implementation initiative for Banking and Banking Operations to achieve Digital transformation through modern banking operations and transaction processing. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization."

# Project overview block
Replace-All "Project Name: Finance and Machine Learning Implementation Initiative" "Project Name: Banking and Banking Operations Implementation Initiative"
Replace-All "Industry Focus: Finance and Machine Learning" "Industry Focus: Banking and Banking Operations"
Replace-All "Digital transformation through intelligent automation and predictive analytics" "Digital transformation through modern banking operations and transaction processing"
Replace-All "This project directly supports organizational strategic objectives by addressing critical business challenges and enabling competitive differentiation through Finance Implementation capabilities." "This project directly supports organizational strategic objectives by addressing critical business challenges and enabling competitive differentiation through Banking Implementation capabilities."

# Business case - proposed solution
Replace-All "Implementation of comprehensive Finance Implementation solution leveraging industry-leading technologies and best practices to address current challenges and enable future growth." "Implementation of comprehensive Banking Implementation solution leveraging industry-leading technologies and best practices to address current challenges and enable future growth."

# Add a (blank) default header and footer to the document's only section.
$section = $d.Sections(1)
$header = $section.Headers(1)
$header.Range.Style = "Header"
$footer = $section.Footers(1)
$footer.Range.Style = "Footer"
